$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Range("E$row")
    if ($cell.Value2 -eq "fullRNASEQ") {
        $cell.Value2 = "fullRNASeq"
    }
}
